# Replace OIE -> WOAH across the workbook's text content (commit message:
# "OIE replaced with WOAH all Excels").
#
# Every occurrence below was located by searching the shared-string table for
# the literal substring "OIE" and mapping each shared string back to the
# single cell that references it.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("References")

$ws1.Range("E5").Value = "Based on official disease reports to the WOAH"
$ws1.Range("E6").Value = "Nairobi sheep disease is listed in the World Organisation for Animal Health ({ref005:WOAH}) Terrestrial Animal Health Code. The map to the right displays occurrence reported to the {ref001:WOAH-WAHIS} system since 2005."
$ws1.Range("E7").Value = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"
$ws1.Range("E14").Value = "Information on stable situations (disease present or absent in a zone or country) is provided by countries through the WOAH monitoring system, which is a different reporting channel. This information is available in a different spatial and temporal scale, which can be browsed on the map independently from the outbreak notification points."
$ws1.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."
$ws1.Range("E72").Value = "Diagnosis of Bunyaviral diseases are described in a dedicated chapter of the {ref010:WOAH Terrestrial Manual}."
$ws1.Range("E137").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

$ws2.Range("C2").Value = "WOAH-WAHIS (WOAH World Animal Health Information System)"
$ws2.Range("C5").Value = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$ws2.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2021. Bunyaviral diseases. Chapter 3.7.6. WOAH Terrestrial Manual, Paris, France"
